$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 7
$ws.Range("H7").Value = 50000
$ws.Range("I7").Value = 50000
$ws.Range("K7").Value = 50000
$ws.Range("M7").Value = -49888
# Row 14
$ws.Range("H14").Value = 50000
$ws.Range("I14").Value = 50000
$ws.Range("K14").Value = 50000
$ws.Range("M14").Value = -49809
# Row 21
$ws.Range("H21").Value = 15000
$ws.Range("J21").Value = 15000
$ws.Range("L21").Value = 15000
$ws.Range("N21").Value = -15936
# Row 23
$ws.Range("H23").Value = 15000
$ws.Range("J23").Value = 15000
$ws.Range("L23").Value = 15000
$ws.Range("N23").Value = -15468
# Row 29
$ws.Range("H29").Value = 7474
$ws.Range("I29").Value = 5655.5
$ws.Range("J29").Value = 11111
$ws.Range("K29").Value = 16966.5
$ws.Range("L29").Value = 33333
$ws.Range("M29").Value = -16685.5
$ws.Range("N29").Value = -33895
# Row 38
$ws.Range("H38").Value = 306.16666
$ws.Range("I38").Value = 37.4
$ws.Range("J38").Value = 498.14285
$ws.Range("K38").Value = 112.2
$ws.Range("L38").Value = 1494.42855
$ws.Range("M38").Value = 259.8
$ws.Range("N38").Value = -2238.42855
# Row 76
$ws.Range("H76").Value = 3370254.8
$ws.Range("I76").Value = 4118319
$ws.Range("J76").Value = 3966.6667
$ws.Range("K76").Value = 4118319
$ws.Range("L76").Value = 3966.6667
$ws.Range("M76").Value = -4118004
$ws.Range("N76").Value = -4596.6667
# Row 79
$ws.Range("H79").Value = 3370254.8
$ws.Range("I79").Value = 4118319
$ws.Range("J79").Value = 3966.6667
$ws.Range("K79").Value = 4118319
$ws.Range("L79").Value = 3966.6667
$ws.Range("M79").Value = -4117227
$ws.Range("N79").Value = -6150.6667
# Row 112
$ws.Range("H112").Value = 7178259.5
$ws.Range("J112").Value = 7793479
$ws.Range("L112").Value = 23380437
$ws.Range("N112").Value = -23382653
# Row 138
$ws.Range("H138").Value = 13633150
$ws.Range("J138").Value = 25003654
$ws.Range("L138").Value = 75010962
$ws.Range("N138").Value = -75021242

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 34399.035
$ws.Range("I2").Value = 35393.656
$ws.Range("J2").Value = 5555
$ws.Range("K2").Value = 35393.656
$ws.Range("L2").Value = 5555
$ws.Range("M2").Value = -35280.656
$ws.Range("N2").Value = -5781
# Row 18
$ws.Range("H18").Value = 50000
$ws.Range("J18").Value = 50000
$ws.Range("L18").Value = 50000
$ws.Range("N18").Value = -50644
# Row 32
$ws.Range("H32").Value = 17333.518
$ws.Range("I32").Value = 3683.043
$ws.Range("J32").Value = 67624.734
$ws.Range("K32").Value = 3683.043
$ws.Range("L32").Value = 67624.734
$ws.Range("M32").Value = -3396.043
$ws.Range("N32").Value = -68198.734
# Row 74
$ws.Range("H74").Value = 4978.0586
$ws.Range("I74").Value = 960.6
$ws.Range("K74").Value = 960.6
$ws.Range("M74").Value = -86.60000000000002
# Row 77
$ws.Range("H77").Value = 4978.0586
$ws.Range("I77").Value = 960.6
$ws.Range("K77").Value = 4803
$ws.Range("M77").Value = -435
# Row 116
$ws.Range("H116").Value = 34399.035
$ws.Range("I116").Value = 35393.656
$ws.Range("J116").Value = 5555
$ws.Range("K116").Value = 35393.656
$ws.Range("L116").Value = 5555
$ws.Range("M116").Value = -33099.656
$ws.Range("N116").Value = -10143

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 34399.035
$ws.Range("I3").Value = 35393.656
$ws.Range("J3").Value = 5555
$ws.Range("K3").Value = 35393.656
$ws.Range("L3").Value = 5555
$ws.Range("M3").Value = -35279.656
$ws.Range("N3").Value = -5783

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6401.3335
$ws.Range("I31").Value = 2073.56
$ws.Range("J31").Value = 14129.5
$ws.Range("K31").Value = 2073.56
$ws.Range("L31").Value = 14129.5
$ws.Range("M31").Value = -1778.56
$ws.Range("N31").Value = -14719.5
# Row 34
$ws.Range("H34").Value = 6401.3335
$ws.Range("I34").Value = 2073.56
$ws.Range("J34").Value = 14129.5
$ws.Range("K34").Value = 2073.56
$ws.Range("L34").Value = 14129.5
$ws.Range("M34").Value = -1871.56
$ws.Range("N34").Value = -14533.5
# Row 107
$ws.Range("H107").Value = 526.55554
$ws.Range("I107").Value = 447.6
$ws.Range("J107").Value = 625.25
$ws.Range("K107").Value = 447.6
$ws.Range("L107").Value = 625.25
$ws.Range("M107").Value = 1472.4
$ws.Range("N107").Value = -4465.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 226.28572
$ws.Range("I6").Value = 97.333336
$ws.Range("J6").Value = 1000
$ws.Range("K6").Value = 292.000008
$ws.Range("L6").Value = 3000
$ws.Range("M6").Value = -179.000008
$ws.Range("N6").Value = -3226
# Row 23
$ws.Range("H23").Value = 842.4
$ws.Range("J23").Value = 313.83334
$ws.Range("L23").Value = 941.5000200000001
$ws.Range("N23").Value = -1411.50002
# Row 113
$ws.Range("H113").Value = 560.0952
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 560.0952
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1680.2856
$ws.Range("M113").Value = ""
$ws.Range("N113").Value = -6020.2856
# Row 114
$ws.Range("H114").Value = 776.2857
$ws.Range("I114").Value = 186.8
$ws.Range("J114").Value = 2250
$ws.Range("K114").Value = 560.4000000000001
$ws.Range("L114").Value = 6750
$ws.Range("M114").Value = 2693.6
$ws.Range("N114").Value = -13258
# Row 117
$ws.Range("H117").Value = 600
$ws.Range("I117").Value = 600
$ws.Range("K117").Value = 1800
$ws.Range("M117").Value = 1642
# Row 131
$ws.Range("H131").Value = 7753550.5
$ws.Range("J131").Value = 8131746
$ws.Range("L131").Value = 24395238
$ws.Range("N131").Value = -24405318

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 40
$ws.Range("H40").Value = 6734
$ws.Range("I40").Value = 5016
$ws.Range("J40").Value = 6948.75
$ws.Range("K40").Value = 5016
$ws.Range("L40").Value = 6948.75
$ws.Range("M40").Value = -4865
$ws.Range("N40").Value = -7250.75
# Row 46
$ws.Range("H46").Value = 17822.223
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 19925
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 19925
$ws.Range("M46").Value = -844
$ws.Range("N46").Value = -20237

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 7534.375
$ws.Range("I22").Value = 1041.6666
$ws.Range("J22").Value = 11430
$ws.Range("K22").Value = 1041.6666
$ws.Range("L22").Value = 11430
$ws.Range("M22").Value = -746.6666
$ws.Range("N22").Value = -12020
# Row 27
$ws.Range("H27").Value = 7534.375
$ws.Range("I27").Value = 1041.6666
$ws.Range("J27").Value = 11430
$ws.Range("K27").Value = 1041.6666
$ws.Range("L27").Value = 11430
$ws.Range("M27").Value = -934.6666
$ws.Range("N27").Value = -11644
# Row 46
$ws.Range("H46").Value = 1962.375
$ws.Range("I46").Value = 1499.75
$ws.Range("K46").Value = 1499.75
$ws.Range("M46").Value = -1311.75
# Row 61
$ws.Range("H61").Value = 1534.7878
$ws.Range("I61").Value = 1297.625
$ws.Range("J61").Value = 2167.2222
$ws.Range("K61").Value = 1297.625
$ws.Range("L61").Value = 2167.2222
$ws.Range("M61").Value = -1095.625
$ws.Range("N61").Value = -2571.2222
# Row 104
$ws.Range("H104").Value = 25400
$ws.Range("J104").Value = 25400
$ws.Range("L104").Value = 25400
$ws.Range("N104").Value = -32388
# Row 113
$ws.Range("H113").Value = 1534.7878
$ws.Range("I113").Value = 1297.625
$ws.Range("J113").Value = 2167.2222
$ws.Range("K113").Value = 1297.625
$ws.Range("L113").Value = 2167.2222
$ws.Range("M113").Value = 872.375
$ws.Range("N113").Value = -6507.2222
# Row 132
$ws.Range("H132").Value = 2926.578
$ws.Range("I132").Value = 1787.2667
$ws.Range("J132").Value = 5205.2
$ws.Range("K132").Value = 5361.800099999999
$ws.Range("L132").Value = 15615.6
$ws.Range("M132").Value = -2831.800099999999
$ws.Range("N132").Value = -20675.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 3823.923
$ws.Range("I81").Value = 2183.3333
$ws.Range("J81").Value = 4692.4707
$ws.Range("K81").Value = 4366.6666
$ws.Range("L81").Value = 9384.9414
$ws.Range("M81").Value = -3305.6666
$ws.Range("N81").Value = -11506.9414
# Row 84
$ws.Range("H84").Value = 3823.923
$ws.Range("I84").Value = 2183.3333
$ws.Range("J84").Value = 4692.4707
$ws.Range("K84").Value = 21833.333
$ws.Range("L84").Value = 46924.70699999999
$ws.Range("M84").Value = -16529.333
$ws.Range("N84").Value = -57532.70699999999
# Row 136
$ws.Range("H136").Value = 1541.96
$ws.Range("I136").Value = 916.619
$ws.Range("J136").Value = 2749.857
$ws.Range("L136").Value = 14475
$ws.Range("M136").Value = -199.857
$ws.Range("N136").Value = -19575
